$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.157153
$ws.Range("H2").Value = 3.471459
$ws.Range("I2").Value = 0.7933642744864839
$ws.Range("J2").Value = 0.7933642744864839
$ws.Range("M2").Value = 30.52246933333333
$ws.Range("N2").Value = 91.567408
$ws.Range("O2").Value = 0.1058764512547768
$ws.Range("P2").Value = 0.1058764512547769
$ws.Range("Q2").Value = 35.31916695647467
$ws.Range("R2").Value = 317.872502608272
$ws.Range("S2").Value = 0.08399859393494961
$ws.Range("T2").Value = 0.08399859393494963
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.157153
$ws.Range("H3").Value = 3.471459
$ws.Range("I3").Value = 0.7933642744864839
$ws.Range("J3").Value = 0.7933642744864839
$ws.Range("O3").Value = 0.001067503492562006
$ws.Range("P3").Value = 0.001067503492562006
$ws.Range("Q3").Value = 0.356106892832
$ws.Range("R3").Value = 3.204962035488
$ws.Range("S3").Value = 0.0008469191338882433
$ws.Range("T3").Value = 0.0008469191338882434
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.157153
$ws.Range("H4").Value = 3.471459
$ws.Range("I4").Value = 0.7933642744864839
$ws.Range("J4").Value = 0.7933642744864839
$ws.Range("M4").Value = 47.57542166666667
$ws.Range("N4").Value = 142.726265
$ws.Range("O4").Value = 0.1650297935598315
$ws.Range("P4").Value = 0.1650297935598315
$ws.Range("Q4").Value = 55.05204190784834
$ws.Range("R4").Value = 495.4683771706351
$ws.Range("S4").Value = 0.13092874243625
$ws.Range("T4").Value = 0.13092874243625
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.157153
$ws.Range("H5").Value = 3.471459
$ws.Range("I5").Value = 0.7933642744864839
$ws.Range("J5").Value = 0.7933642744864839
$ws.Range("M5").Value = 209.8781993333333
$ws.Range("N5").Value = 629.634598
$ws.Range("O5").Value = 0.7280262516928295
$ws.Range("P5").Value = 0.7280262516928296
$ws.Range("Q5").Value = 242.8611879931647
$ws.Range("R5").Value = 2185.750691938482
$ws.Range("S5").Value = 0.5775900189813961
$ws.Range("T5").Value = 0.5775900189813962
$ws.Range("G6").Value = 0.3013863333333334
$ws.Range("H6").Value = 0.904159
$ws.Range("I6").Value = 0.2066357255135161
$ws.Range("J6").Value = 0.206635725513516
$ws.Range("M6").Value = 30.52246933333333
$ws.Range("N6").Value = 91.567408
$ws.Range("O6").Value = 0.1058764512547768
$ws.Range("P6").Value = 0.1058764512547769
$ws.Range("Q6").Value = 9.199055116652445
$ws.Range("R6").Value = 82.79149604987201
$ws.Range("S6").Value = 0.02187785731982723
$ws.Range("T6").Value = 0.02187785731982723
$ws.Range("G7").Value = 0.3013863333333334
$ws.Range("H7").Value = 0.904159
$ws.Range("I7").Value = 0.2066357255135161
$ws.Range("J7").Value = 0.206635725513516
$ws.Range("O7").Value = 0.001067503492562006
$ws.Range("P7").Value = 0.001067503492562006
$ws.Range("Q7").Value = 0.09274983576533333
$ws.Range("R7").Value = 0.834748521888
$ws.Range("S7").Value = 0.0002205843586737623
$ws.Range("T7").Value = 0.0002205843586737623
$ws.Range("G8").Value = 0.3013863333333334
$ws.Range("H8").Value = 0.904159
$ws.Range("I8").Value = 0.2066357255135161
$ws.Range("J8").Value = 0.206635725513516
$ws.Range("M8").Value = 47.57542166666667
$ws.Range("N8").Value = 142.726265
$ws.Range("O8").Value = 0.1650297935598315
$ws.Range("P8").Value = 0.1650297935598315
$ws.Range("Q8").Value = 14.33858189290389
$ws.Range("R8").Value = 129.047237036135
$ws.Range("S8").Value = 0.03410105112358157
$ws.Range("T8").Value = 0.03410105112358156
$ws.Range("G9").Value = 0.3013863333333334
$ws.Range("H9").Value = 0.904159
$ws.Range("I9").Value = 0.2066357255135161
$ws.Range("J9").Value = 0.206635725513516
$ws.Range("M9").Value = 209.8781993333333
$ws.Range("N9").Value = 629.634598
$ws.Range("O9").Value = 0.7280262516928295
$ws.Range("P9").Value = 0.7280262516928296
$ws.Range("Q9").Value = 63.25442094367578
$ws.Range("R9").Value = 569.289788493082
$ws.Range("S9").Value = 0.1504362327114335
$ws.Range("T9").Value = 0.1504362327114335
